# Evaluations: #17 Part 1. CSS selectors by functions in AST.
#
# - "Errata" sheet (Sheet2) gets a new errata row (Entry 3 / "title").
# - "Selectors" sheet (Sheet1) gets a new "1 Found" column inserted just
#   before the existing "3 Found" column, flagged "y" for every data row.
# - The active sheet/tab switches from "Errata" back to "Selectors", and
#   the "Errata" sheet's remembered selection moves to A5.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Selectors"
$ws2 = $wb.Worksheets.Item(2)   # "Errata"

# --- Errata sheet: append row 4 (do this first so the shared-string
#     table picks up "title" before "1 Found", matching insertion order) ---
$ws2.Cells.Item(4, 1).Value = 3
$ws2.Cells.Item(4, 2).Value = "title"

# --- Selectors sheet: insert a new "1 Found" column before column G ---
$ws1.Columns.Item(7).EntireColumn.Insert()
$ws1.Columns.Item(7).ColumnWidth = 8

$ws1.Cells.Item(1, 7).Value = "1 Found"
for ($r = 2; $r -le 17; $r++) {
    $ws1.Cells.Item($r, 7).Value = "y"
}

# --- View state: Errata's selection moves to A5, then Selectors becomes
#     the active/selected tab again ---
$ws2.Activate() | Out-Null
$ws2.Range("A5").Select() | Out-Null
$ws1.Activate() | Out-Null
